# Apply commit "add substr and start part two":
#  - Append a new "Second Part" section below the existing "First part" table
#    (mirrors the structure of the existing sections, e.g. rows 7/35 headers
#    and rows 8-29/36-37 data rows), with the first new function (ft_substr)
#    already implemented, and the rest of part two started but not done.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 42: section title "Second Part" (same plain style as "First part" in D5) ---
$ws.Range("D42").Value = "Second Part"

# --- Row 44: table header, copy the look of the existing header row (D35:I35) ---
$ws.Range("D35:I35").Copy()
$ws.Range("D44:I44").PasteSpecial(-4122)
$ws.Range("B44").Value = 1
$ws.Range("D44").Value = "Mandatory functions"
$ws.Range("E44").Value = "Implementation"
$ws.Range("F44").Value = "Test"
$ws.Range("G44").Value = "Norminette"
$ws.Range("H44").Value = "Done Flag"
$ws.Range("I44").Value = "Comment"

# --- Rows 45-54: data rows, copy the look of an existing data row (D36:I36) ---
$ws.Range("D36:I36").Copy()
$ws.Range("D45:I54").PasteSpecial(-4122)

for ($r = 45; $r -le 54; $r++) {
    $ws.Range("B$r").Value = 1
}

$ws.Range("D45").Value = "ft_substr"
$ws.Range("E45").Value = "Done"
$ws.Range("F45").Value = "Valid"
$ws.Range("G45").Value = "Not pass"
$ws.Range("H45").Value = 0.75

$ws.Range("D46").Value = "ft_strjoin"
$ws.Range("E46").Value = "Not Done"
$ws.Range("F46").Value = "Not Valid"
$ws.Range("G46").Value = "Not pass"
$ws.Range("H46").Value = 0

$ws.Range("D47").Value = "ft_strtrim"
$ws.Range("E47").Value = "Not Done"
$ws.Range("F47").Value = "Not Valid"
$ws.Range("G47").Value = "Not pass"
$ws.Range("H47").Value = 0

$ws.Range("D48").Value = "ft_split"
$ws.Range("E48").Value = "Not Done"
$ws.Range("F48").Value = "Not Valid"
$ws.Range("G48").Value = "Not pass"
$ws.Range("H48").Value = 0

$ws.Range("D49").Value = "ft_itoa"
$ws.Range("E49").Value = "Not Done"
$ws.Range("F49").Value = "Not Valid"
$ws.Range("G49").Value = "Not pass"
$ws.Range("H49").Value = 0

$ws.Range("D50").Value = "ft_strmapi"
$ws.Range("E50").Value = "Not Done"
$ws.Range("F50").Value = "Not Valid"
$ws.Range("G50").Value = "Not pass"
$ws.Range("H50").Value = 0

$ws.Range("D51").Value = "ft_putchar_fd"
$ws.Range("E51").Value = "Not Done"
$ws.Range("F51").Value = "Not Valid"
$ws.Range("G51").Value = "Not pass"
$ws.Range("H51").Value = 0

$ws.Range("D52").Value = "ft_putstr_fd"
$ws.Range("E52").Value = "Not Done"
$ws.Range("F52").Value = "Not Valid"
$ws.Range("G52").Value = "Not pass"
$ws.Range("H52").Value = 0

$ws.Range("D53").Value = "ft_putendl_fd"
$ws.Range("E53").Value = "Not Done"
$ws.Range("F53").Value = "Not Valid"
$ws.Range("G53").Value = "Not Pass"
$ws.Range("H53").Value = 0

$ws.Range("D54").Value = "ft_putnbr_fd"
$ws.Range("E54").Value = "Not Done"
$ws.Range("F54").Value = "Not Valid"
$ws.Range("G54").Value = "Not Pass"
$ws.Range("H54").Value = 0

# --- Row 55: totals for the new section (same pattern as row 30) ---
$ws.Range("B55").Formula = "=SUM(B44:B54)"
$ws.Range("H55").Formula = "=SUM(H45:H54)"

# --- Row 57-58: percentage of completion for the new section (same pattern as rows 32-33) ---
$ws.Range("I32:I33").Copy()
$ws.Range("I57:I58").PasteSpecial(-4122)
$ws.Range("I57").Value = "Percentage of completion"
$ws.Range("I58").Formula = "=H55/B55"

# --- Extend the "na" conditional formatting (dxfId=2) to also cover the new data block ---
$newNaRule = $ws.Range("E45:G54").FormatConditions.Add(1, 3, """na""")
$newNaRule.Font.Color = 22428
$newNaRule.Interior.Color = 10284031

# --- Update the active selection / view to match where the edit ended up ---
$ws.Range("H56").Select()
